$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values per 28C discharge temperature sizing recalculation ---
$ws.Range("B11").Value = 2.2599999999999998
$ws.Range("C11").Value = 0.33142100000000002
$ws.Range("D11").Value = 0.33900000000000002
$ws.Range("E11").Value = 0.76659500000000003
$ws.Range("F11").Value = 0.00077899999999999996
$ws.Range("G11").Value = 0.33900000000000002
$ws.Range("B12").Value = 1.1000000000000001
$ws.Range("C12").Value = 0.14306099999999999
$ws.Range("D12").Value = 0.30399999999999999
$ws.Range("E12").Value = 0.33572800000000003
$ws.Range("F12").Value = 0.001619
$ws.Range("G12").Value = 0.30399999999999999
$ws.Range("B13").Value = 0.789493
$ws.Range("C13").Value = 0.126691
$ws.Range("D13").Value = 0.2
$ws.Range("E13").Value = 0.15789900000000001
$ws.Range("F13").Value = 0.00120300000000000008
$ws.Range("G13").Value = 0.2
$ws.Range("B14").Value = 1.1200000000000001
$ws.Range("C14").Value = 0.14118700000000001
$ws.Range("D14").Value = 0.249
$ws.Range("E14").Value = 0.27863700000000002
$ws.Range("F14").Value = 0.00134399999999999993
$ws.Range("G14").Value = 0.249
$ws.Range("B15").Value = 0.790601
$ws.Range("C15").Value = 0.12650500000000001
$ws.Range("D15").Value = 0.2
$ws.Range("E15").Value = 0.15812000000000001
$ws.Range("F15").Value = 0.00120499999999999991
$ws.Range("G15").Value = 0.2
$ws.Range("B22").Value = 31889.35
$ws.Range("B23").Value = 20127.95
$ws.Range("B24").Value = 13248.64
$ws.Range("B25").Value = 20395.03
$ws.Range("B26").Value = 13414.6
$ws.Range("B33").Value = 1.7
$ws.Range("C33").Value = 1.7
$ws.Range("D33").Value = 6.06
$ws.Range("E33").Value = 5.44
$ws.Range("F33").Value = 4.5199999999999996
$ws.Range("G33").Value = 5.44
$ws.Range("I33").Value = 1
$ws.Range("J33").Value = 5.44
$ws.Range("B39").Value = 5.44
$ws.Range("B44").Value = 5.44
$ws.Range("C44").Value = 128026.53
$ws.Range("D44").Value = 0.71741200000000005
$ws.Range("E44").Value = 1.81
$ws.Range("F44").Value = 42671.24
$ws.Range("G44").Value = 0.71741200000000005
$ws.Range("B50").Value = 43895.519999999997
$ws.Range("B56").Value = 5.44
$ws.Range("C69").Value = 21.84
$ws.Range("C70").Value = 4.6100000000000003
$ws.Range("C71").Value = 2.56
$ws.Range("C72").Value = 4.6100000000000003
$ws.Range("C73").Value = 2.92
$ws.Range("D80").Value = 2.2612999999999999
$ws.Range("E80").Value = 2.2612999999999999
$ws.Range("F80").Value = 0.76659999999999995
$ws.Range("D81").Value = 1.1285000000000001
$ws.Range("E81").Value = 1.1285000000000001
$ws.Range("F81").Value = 0.40629999999999999
$ws.Range("D82").Value = 0.74139999999999995
$ws.Range("E82").Value = 0.74139999999999995
$ws.Range("F82").Value = 0.14829999999999999
$ws.Range("D83").Value = 1.139
$ws.Range("E83").Value = 1.139
$ws.Range("F83").Value = 0.35199999999999998
$ws.Range("D84").Value = 0.75009999999999999
$ws.Range("E84").Value = 0.75009999999999999
$ws.Range("F84").Value = 0.15
$ws.Range("B90").Value = 6.0648
$ws.Range("C90").Value = 4.5202
$ws.Range("D90").Value = 6.0648
$ws.Range("E90").Value = 1.6970000000000001
$ws.Range("D96").Value = 0.76659999999999995
$ws.Range("E96").Value = 0.76659999999999995
$ws.Range("F96").Value = 0.76659999999999995
$ws.Range("D97").Value = 0.3357
$ws.Range("E97").Value = 0.3357
$ws.Range("F97").Value = 0.3357
$ws.Range("D98").Value = 0.15790000000000001
$ws.Range("E98").Value = 0.15790000000000001
$ws.Range("F98").Value = 0.15790000000000001
$ws.Range("D99").Value = 0.27860000000000001
$ws.Range("E99").Value = 0.27860000000000001
$ws.Range("F99").Value = 0.27860000000000001
$ws.Range("D100").Value = 0.15809999999999999
$ws.Range("E100").Value = 0.15809999999999999
$ws.Range("F100").Value = 0.15809999999999999
$ws.Range("B106").Value = 1.6970000000000001
$ws.Range("C106").Value = 1.6970000000000001
$ws.Range("D106").Value = 1.6970000000000001
$ws.Range("E106").Value = 1.6970000000000001

# --- Row height updates (rows 23-26 grew from 36 to 54 to fit new wrapped text) ---
$ws.Rows.Item(23).RowHeight = 54
$ws.Rows.Item(24).RowHeight = 54
$ws.Rows.Item(25).RowHeight = 54
$ws.Rows.Item(26).RowHeight = 54

# --- Sheet view: zoom out to 40% and move selection to F85 ---
$ws.Select()
$excel.ActiveWindow.Zoom = 40
$ws.Range("F85").Select()
